$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 129. This shifts the existing rows
# 129 -> 130 and 130 -> 131 (their values/styles travel with them), and the
# new row 129 inherits formatting (e.g. the date style on column D) from the
# row above, matching Excel's normal insert behaviour.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new record.
$ws.Cells.Item(129, 1).Value = 10
$ws.Cells.Item(129, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(129, 3).Value = "La Araucanía"
$ws.Cells.Item(129, 4).Value2 = [DateTime]"2021-09-09"
$ws.Cells.Item(129, 5).Value = 9
$ws.Cells.Item(129, 6).Value = "Fruta"
$ws.Cells.Item(129, 7).Value = 100102
$ws.Cells.Item(129, 8).Value = "Cítricos"
$ws.Cells.Item(129, 9).Value = 100102006
$ws.Cells.Item(129, 10).Value = "Pomelo"
$ws.Cells.Item(129, 11).Value = "Start Ruby"
$ws.Cells.Item(129, 12).Value = "Primera"
$ws.Cells.Item(129, 13).Value = 240
$ws.Cells.Item(129, 14).Value = 12000
$ws.Cells.Item(129, 15).Value = 13000
$ws.Cells.Item(129, 16).Value = 12583
$ws.Cells.Item(129, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(129, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(129, 19).Value = 839
$ws.Cells.Item(129, 20).Value = 15
